$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2,1).Value = "Última actualización: 02:49:45"
$ws.Cells.Item(3,1).Value = "Total filas: 10"
$ws.Cells.Item(8,1).Value = "02:49:45"
$ws.Cells.Item(8,2).Value = "02:49"
$ws.Cells.Item(8,4).Value = 0
$ws.Cells.Item(9,1).Value = "01:12:01"
$ws.Cells.Item(9,2).Value = "02:58"
$ws.Cells.Item(9,4).Value = 106
$ws.Cells.Item(10,2).Value = "02:59"
$ws.Cells.Item(10,3).Value = "215_ALUAR"
$ws.Cells.Item(10,4).Value = 63
$ws.Cells.Item(11,1).Value = "02:49:45"
$ws.Cells.Item(11,2).Value = "03:48"
$ws.Cells.Item(11,4).Value = 59
$ws.Cells.Item(12,2).Value = "03:53"
$ws.Cells.Item(12,3).Value = "14_ABASTO"
$ws.Cells.Item(12,4).Value = 89
$ws.Cells.Item(13,2).Value = "03:58"
$ws.Cells.Item(13,3).Value = "215_ALUAR"
$ws.Cells.Item(13,4).Value = 94
$ws.Cells.Item(14,1).Value = "02:49:45"
$ws.Cells.Item(14,2).Value = "04:01"
$ws.Cells.Item(14,3).Value = "81_EL PELIGRO"
$ws.Cells.Item(14,4).Value = 72
$ws.Cells.Item(14,5).Value = "LP1912"
$ws.Cells.Item(15,1).Value = "02:49:45"
$ws.Cells.Item(15,2).Value = "04:35"
$ws.Cells.Item(15,3).Value = "215_ALUAR"
$ws.Cells.Item(15,4).Value = 106
$ws.Cells.Item(15,5).Value = "LP1912"

# --- Sheet: LP1912-215 ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2,1).Value = "Última actualización: 02:49:45"
$ws.Cells.Item(3,1).Value = "Total filas: 6"
$ws.Cells.Item(7,1).Value = "02:49:45"
$ws.Cells.Item(7,2).Value = "02:49"
$ws.Cells.Item(7,4).Value = 0
$ws.Cells.Item(8,1).Value = "01:12:01"
$ws.Cells.Item(8,2).Value = "02:58"
$ws.Cells.Item(8,4).Value = 106
$ws.Cells.Item(9,1).Value = "01:56:31"
$ws.Cells.Item(9,2).Value = "02:59"
$ws.Cells.Item(9,4).Value = 63
$ws.Cells.Item(10,1).Value = "02:24:16"
$ws.Cells.Item(10,2).Value = "03:58"
$ws.Cells.Item(10,3).Value = "215_ALUAR"
$ws.Cells.Item(10,4).Value = 94
$ws.Cells.Item(10,5).Value = "LP1912"
$ws.Cells.Item(11,1).Value = "02:49:45"
$ws.Cells.Item(11,2).Value = "04:35"
$ws.Cells.Item(11,3).Value = "215_ALUAR"
$ws.Cells.Item(11,4).Value = 106
$ws.Cells.Item(11,5).Value = "LP1912"

# --- Sheet: 6203-6173 ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2,1).Value = "Última actualización: 02:49:45"
